# Weekly update: a new price-observation row is inserted at row 248 of the
# "Hortaliza, Terminal Hortofrutícola Agro Chillán - Zanahoria" sheet. All
# rows from the former 248 onward shift down by one (248 -> 249, ..., 314 ->
# 315) and the sheet's used range grows from A1:R314 to A1:R315.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 248 - this pushes the existing row 248
# (and everything below it) down to row 249, etc.
$ws.Rows("248").Insert()

# Populate the newly inserted row 248 with the new observation.
$ws.Range("A248").Value = 7
$ws.Range("B248").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C248").Value = "Ñuble"
$ws.Range("D248").Value = (Get-Date -Year 2022 -Month 8 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E248").Value = 16
$ws.Range("F248").Value = 100114013
$ws.Range("G248").Value = "Zanahoria"
$ws.Range("H248").Value = "Sin especificar"
$ws.Range("I248").Value = "Primera"
$ws.Range("J248").Value = 160
$ws.Range("K248").Value = 7500
$ws.Range("L248").Value = 8000
$ws.Range("M248").Value = 7750
$ws.Range("N248").Value = '$/saco 20 kilos'
$ws.Range("O248").Value = "Región de Ñuble"
$ws.Range("P248").Value = 388
$ws.Range("Q248").Value = 20
$ws.Range("R248").Value = "Hortaliza"
